$wb = $excel.ActiveWorkbook

# "Computational benchmarks" sheet holds the Sphere / ITER 1D benchmark
# rows that are being switched on for post-processing.
$ws = $wb.Worksheets.Item("Computational benchmarks")
$ws.Activate()

# The "Run" / "OnlyInput" columns store their true/false flags as literal
# text (not native booleans), so instead of assigning the strings "true"/
# "false" directly (which the engine auto-coerces into boolean cells), copy
# them from existing cells that already hold that literal text. This keeps
# the cell type/style identical to the rest of the sheet.

# Row 4 (Sphere Leakage Test): turn "Run" (D4) and "OnlyInput" (E4) on.
# E5 already contains the literal text "true" - use it as the source.
$ws.Range("E5").Copy()
$ws.Range("D4").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E4").PasteSpecial(-4163)  # xlPasteValues

# Row 5 (ITER 1D): turn "OnlyInput" (E5) off, using a cell that already
# contains the literal text "false".
$ws.Range("C4").Copy()
$ws.Range("E5").PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = $false

# Move the active selection from D6 to E6, beginning the post-processing work.
$ws.Range("E6").Select()
